$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography ("o século 21. Editora
# Campus. 2010."). Immediately after it there are two empty paragraphs and
# then the "© 2020 ..." footer paragraph, all of which must be removed.
$r = $d.Content
$found = $r.Find.Execute("o século 21. Editora Campus. 2010.", $true,
                          $false, $false, $false, $false, $true, 1,
                          $false, "", 0)

if ($found) {
    $anchorPara = $r.Paragraphs(1)

    $emptyPara1 = $anchorPara.Next()
    $emptyPara2 = $emptyPara1.Next()
    $copyrightPara = $emptyPara2.Next()

    $deleteRange = $d.Range($emptyPara1.Range.Start, $copyrightPara.Range.End)
    $deleteRange.Delete()
}
